$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "60.296.43"
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.333.31"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  +0.00%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "547.98"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "130.99"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.21%  "
$ws.Range("E7").Value = "  +0.04%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.580"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.59%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.330.03"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("E10").Value = "  +0.57%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.63"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("E12").Value = "  -0.68%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.337"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.38%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "23.70"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.96%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "2.747.20"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.58%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "60.258.25"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("E17").Value = "  +0.81%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.332.88"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.69%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "10.69"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("E20").Value = "  -1.68%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "315.12"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.16%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.59"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -3.50%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.11%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "64.19"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("E25").Value = "  -0.79%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("E28").Value = "  +0.96%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.27"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +9.82%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.73"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "171.54"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.09%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.0₃0733"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.24%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "6.07"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("E34").Value = "  -3.27%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.384"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.28%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "18.08"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").Value = "  +0.00%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.06%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "4.11"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.81%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "322.70"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.91%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "38.12"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("E42").Value = "  -0.64%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "137.75"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("E44").Value = "  +1.34%  "
$ws.Range("E45").Value = "  -0.12%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "19.36"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("E47").Value = "  +0.71%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0499"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("E49").Value = "  +0.78%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0₆0218"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.76%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "10.94"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.80%  "
